# Reorganise the "Arborescence des pages" tree:
#  - "Menu/réservation" and "Menu/historique" become direct children of
#    "Menu" (new paragraphs inserted right before "Menu/Administrateur").
#  - the old "Menu/réservation" paragraph (the one that used to sit right
#    before "Menu/Administrateur/créationPlace") is removed.
#  - the old "Menu/historique" / "Menu/historiqueComplet" pair (the ones
#    that used to sit right after "Menu/Administrateur/listeAttribution")
#    collapse into a single "Menu/Administrateur/historiqueComplet"
#    paragraph.

$d = $word.ActiveDocument
$cr = [char]13

function Get-ParaIndexByText($doc, [string]$text) {
    # first match, scanning from the top
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text.TrimEnd($cr)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

function Get-LastParaIndexByText($doc, [string]$text) {
    # last match, scanning from the bottom (used once duplicate text
    # exists in the document after earlier insertions)
    $n = $doc.Paragraphs.Count
    for ($i = $n; $i -ge 1; $i--) {
        $t = $doc.Paragraphs.Item($i).Range.Text.TrimEnd($cr)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

$nbsp = [char]0x00A0
$eacute = [char]0x00E9

# ---------------------------------------------------------------------
# 1) Insert "Menu/réservation" and "Menu/historique" as new paragraphs
#    directly before "Menu/Administrateur".
# ---------------------------------------------------------------------
$adminText = "Page" + $nbsp + ": parking.com/Menu/Administrateur"
$adminIdx = Get-ParaIndexByText $d $adminText

$admin = $d.Paragraphs.Item($adminIdx)
$admin.Range.InsertParagraphBefore()
$admin.Range.InsertParagraphBefore()

$d.Paragraphs.Item($adminIdx).Range.Text = "Page : parking.com/Menu/réservation"
$d.Paragraphs.Item($adminIdx + 1).Range.Text = "Page : parking.com/Menu/historique"

# ---------------------------------------------------------------------
# 2) Remove the old "Menu/réservation" paragraph (the one that used to
#    precede "Menu/Administrateur/créationPlace").
# ---------------------------------------------------------------------
$oldResText = "Page" + $nbsp + ": parking.com/Menu/r" + $eacute + "servation"
$oldResIdx = Get-ParaIndexByText $d $oldResText
$d.Paragraphs.Item($oldResIdx).Range.Delete()

# ---------------------------------------------------------------------
# 3) Collapse the old "Menu/historique" + "Menu/historiqueComplet" pair
#    into a single "Menu/Administrateur/historiqueComplet" paragraph.
#    (use the *last* match for "Menu/historique" since step 1 created
#    an earlier paragraph with the same text)
# ---------------------------------------------------------------------
$oldHistIdx = Get-LastParaIndexByText $d "Page : parking.com/Menu/historique"
$oldHistCompletIdx = $oldHistIdx + 1

$d.Paragraphs.Item($oldHistIdx).Range.Text = "Page : parking.com/Menu/Administrateur/historiqueComplet"
$d.Paragraphs.Item($oldHistCompletIdx).Range.Delete()
